# chore: adapt column header formatting to respective input file names
#
# 1. Rename the header-row labels so the "_old"/"_new" suffixes become the
#    actual format-version identifiers ("_FV2304" / "_FV2310").
# 2. Turn the used range into a native Excel Table ("Table1") so the header
#    row gets AutoFilter + structured-reference columns matching the new
#    names.
# 3. Freeze the header row (pane split under row 1) and keep the selection
#    anchored in the (now frozen) lower pane.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1, every used column) -----------------
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count
$lastRow = $usedRange.Rows.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2304"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2310"
        }
    }
}

# --- 2. Convert the used range into a real Table -----------------------
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, [Type]::Missing, 1)
$tbl.Name = "Table1"
# Inherit the workbook's own default table style rather than forcing a
# different look.
$tbl.TableStyle = "TableStyleMedium9"

# --- 3. Freeze the header row ------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
